# Apply "Net Growth" -> "Net Growth Forest" flow-type renaming on the
# "Flow Type-Group Membership" worksheet (rows 57-71, columns A & B).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Flow Type-Group Membership")
$ws.Activate()

# Column A (Flow Type) values for rows 57-71 cycle through these five labels.
$flowTypes = @(
    "Net Growth Forest: Atmosphere -> Coarse Roots",
    "Net Growth Forest: Atmosphere -> Fine Roots",
    "Net Growth Forest: Atmosphere -> Foliage",
    "Net Growth Forest: Atmosphere -> Merchantable",
    "Net Growth Forest: Atmosphere -> Other Wood"
)

for ($row = 57; $row -le 71; $row++) {
    $idx = ($row - 57) % 5
    $ws.Cells.Item($row, 1).Value = $flowTypes[$idx]
}

# Column B (Flow Group) for rows 62-66 changes from "Net Growth: Total"
# to "Net Growth Forest: Total". Rows 57-61 and 67-71 are unchanged.
for ($row = 62; $row -le 66; $row++) {
    $ws.Cells.Item($row, 2).Value = "Net Growth Forest: Total"
}

# Update the sheet selection to a single cell (previously a multi-cell range).
$ws.Range("B62").Select()
